$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append the new "Total Size" column after Postal
$ws.Range("G1").Value = "Total Size"

# Relabel "City" -> "City/ State" (stays in column E)
$ws.Range("E1").Value = "City/ State"

# Selection collapses back to a single cell instead of the pasted A2:XFD3 block
$ws.Range("A2").Select()
